# Update selection on the existing "paymentOptions" sheet.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B7").Select()

# Insert the new "bookNames" worksheet right after "paymentOptions".
$bookNames = $wb.Worksheets.Add($null, $ws1)
$bookNames.Name = "bookNames"

# Populate the header and the book name values.
$bookNames.Range("A1").Value = "bookName"
$bookNames.Range("A2").Value = "CÁCH NỀN KINH TẾ VẬN HÀNH Niềm tin, sự sụp đổ và những lời tiên tri tự đúng"
$bookNames.Range("A3").Value = "Lời Thú Tội Của Một Sát Thủ Kinh Tế - Bìa Cứng (Tái Bản 2023)"
$bookNames.Range("A4").Value = "Tuyển tập Vũ Trọng Phụng"
$bookNames.Range("A5").Value = "Tuyển Tập Truyện Ngắn Hay Nhất Của Nguyễn Minh Châu"
$bookNames.Range("A6").Value = "Văn Học Trong Nhà Trường: Thơ Nguyễn Khuyến"

# Header formatting.
$bookNames.Range("A1").Font.Name = "Aptos Narrow"

# Data rows formatting: font + vertically centered text.
$bookNames.Range("A2:A6").Font.Name = "Aptos Narrow"
$bookNames.Range("A2:A6").VerticalAlignment = -4108

# Widen column A to fit the long titles.
$bookNames.Columns.Item(1).ColumnWidth = 63.16666666666667

# Select a cell on the new sheet and make it the active sheet/tab.
$bookNames.Range("A23").Select()
$bookNames.Activate()
